$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 2 Brazilian Serie B matches in rows 2-3 (dated 2025-10-14).
# The update replaces rows 2-3 with two FIFA World Cup Qualifiers - Asia matches,
# re-inserts the original two Brazilian Serie B matches (with refreshed odds) as new
# rows 4-5, and shifts/updates the existing FIFA World Cup Qualifiers - Americas match
# (with refreshed odds) down to row 6.

# Insert two blank rows at position 4 to make room for the Brazilian Serie B matches,
# pushing the FIFA World Cup Qualifiers - Americas match from row 4 down to row 6.
$ws.Range("A4:A5").EntireRow.Insert()

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = "FIFA World Cup Qualifiers - Asia"
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "2025-10-14"
$ws.Cells.Item(2,2).ClearFormats()
$ws.Cells.Item(2,3).Value = "14:00:00"
$ws.Cells.Item(2,4).Value = "Qatar"
$ws.Cells.Item(2,5).Value = "UAE"
$ws.Cells.Item(2,6).Value = 2.34
$ws.Cells.Item(2,7).Value = 2.38
$ws.Cells.Item(2,8).Value = 3.65
$ws.Cells.Item(2,9).Value = 3.7
$ws.Cells.Item(2,10).Value = 3.2
$ws.Cells.Item(2,11).Value = 3.35
$ws.Cells.Item(2,12).Value = 1.48
$ws.Cells.Item(2,13).Value = 1.07
$ws.Cells.Item(2,14).Value = 3.05
$ws.Cells.Item(2,15).Value = 1.41
$ws.Cells.Item(2,16).Value = 1.71
$ws.Cells.Item(2,17).Value = 2.26
$ws.Cells.Item(2,18).Value = 1.26
$ws.Cells.Item(2,19).Value = 4.2
$ws.Cells.Item(2,20).Value = 1.86
$ws.Cells.Item(2,21).Value = 1.96
$ws.Cells.Item(2,22).Value = 1.31
$ws.Cells.Item(2,23).Value = 1.55
$ws.Cells.Item(2,24).Value = 11.5
$ws.Cells.Item(2,25).Value = 12
$ws.Cells.Item(2,26).Value = 23
$ws.Cells.Item(2,27).Value = 70
$ws.Cells.Item(2,28).Value = 9
$ws.Cells.Item(2,29).Value = 7.6
$ws.Cells.Item(2,30).Value = 15.5
$ws.Cells.Item(2,31).Value = 48
$ws.Cells.Item(2,32).Value = 14.5
$ws.Cells.Item(2,33).Value = 11.5
$ws.Cells.Item(2,34).Value = 20
$ws.Cells.Item(2,35).Value = 65
$ws.Cells.Item(2,36).Value = 36
$ws.Cells.Item(2,37).Value = 29
$ws.Cells.Item(2,38).Value = 48
$ws.Cells.Item(2,39).Value = 130
$ws.Cells.Item(2,40).Value = 26
$ws.Cells.Item(2,41).Value = 55

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = "FIFA World Cup Qualifiers - Asia"
$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,2).Value = "2025-10-14"
$ws.Cells.Item(3,2).ClearFormats()
$ws.Cells.Item(3,3).Value = "15:45:00"
$ws.Cells.Item(3,4).Value = "Saudi Arabia"
$ws.Cells.Item(3,5).Value = "Iraq"
$ws.Cells.Item(3,6).Value = 1.74
$ws.Cells.Item(3,7).Value = 1.88
$ws.Cells.Item(3,8).Value = 4.7
$ws.Cells.Item(3,9).Value = 7.4
$ws.Cells.Item(3,10).Value = 3.3
$ws.Cells.Item(3,11).Value = 4.2
$ws.Cells.Item(3,12).Value = 1.54
$ws.Cells.Item(3,13).Value = 1.1
$ws.Cells.Item(3,14).Value = 2.8
$ws.Cells.Item(3,15).Value = 1.49
$ws.Cells.Item(3,16).Value = 1.6
$ws.Cells.Item(3,17).Value = 2.12
$ws.Cells.Item(3,18).Value = 1.22
$ws.Cells.Item(3,19).Value = 4.5
$ws.Cells.Item(3,20).Value = 1.95
$ws.Cells.Item(3,21).Value = 1.62
$ws.Cells.Item(3,22).Value = 1.15
$ws.Cells.Item(3,23).Value = 1.93
$ws.Cells.Item(3,24).Value = 12.5
$ws.Cells.Item(3,25).Value = 1000
$ws.Cells.Item(3,26).Value = 1000
$ws.Cells.Item(3,27).Value = 1000
$ws.Cells.Item(3,28).Value = 8
$ws.Cells.Item(3,29).Value = 9.800000000000001
$ws.Cells.Item(3,30).Value = 1000
$ws.Cells.Item(3,31).Value = 1000
$ws.Cells.Item(3,32).Value = 1000
$ws.Cells.Item(3,33).Value = 13
$ws.Cells.Item(3,34).Value = 1000
$ws.Cells.Item(3,35).Value = 1000
$ws.Cells.Item(3,36).Value = 1000
$ws.Cells.Item(3,37).Value = 1000
$ws.Cells.Item(3,38).Value = 1000
$ws.Cells.Item(3,39).Value = 1000
$ws.Cells.Item(3,40).Value = 1000
$ws.Cells.Item(3,41).Value = 1000

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = "Brazilian Serie B"
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = "2025-10-14"
$ws.Cells.Item(4,2).ClearFormats()
$ws.Cells.Item(4,3).Value = "19:30:00"
$ws.Cells.Item(4,4).Value = "Paysandu"
$ws.Cells.Item(4,5).Value = "Remo"
$ws.Cells.Item(4,6).Value = 2.9
$ws.Cells.Item(4,7).Value = 3.1
$ws.Cells.Item(4,8).Value = 2.56
$ws.Cells.Item(4,9).Value = 2.76
$ws.Cells.Item(4,10).Value = 3.3
$ws.Cells.Item(4,11).Value = 3.5
$ws.Cells.Item(4,12).Value = 1.42
$ws.Cells.Item(4,13).Value = 1.09
$ws.Cells.Item(4,14).Value = 2.96
$ws.Cells.Item(4,15).Value = 1.44
$ws.Cells.Item(4,16).Value = 1.68
$ws.Cells.Item(4,17).Value = 2.24
$ws.Cells.Item(4,18).Value = 1.25
$ws.Cells.Item(4,19).Value = 4.3
$ws.Cells.Item(4,20).Value = 1.94
$ws.Cells.Item(4,21).Value = 1.93
$ws.Cells.Item(4,22).Value = 1.57
$ws.Cells.Item(4,23).Value = 1.48
$ws.Cells.Item(4,24).Value = 13
$ws.Cells.Item(4,25).Value = 1000
$ws.Cells.Item(4,26).Value = 19
$ws.Cells.Item(4,27).Value = 1000
$ws.Cells.Item(4,28).Value = 11.5
$ws.Cells.Item(4,29).Value = 9
$ws.Cells.Item(4,30).Value = 14
$ws.Cells.Item(4,31).Value = 1000
$ws.Cells.Item(4,32).Value = 1000
$ws.Cells.Item(4,33).Value = 18
$ws.Cells.Item(4,34).Value = 25
$ws.Cells.Item(4,35).Value = 60
$ws.Cells.Item(4,36).Value = 75
$ws.Cells.Item(4,37).Value = 50
$ws.Cells.Item(4,38).Value = 1000
$ws.Cells.Item(4,39).Value = 160
$ws.Cells.Item(4,40).Value = 1000
$ws.Cells.Item(4,41).Value = 50

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = "Brazilian Serie B"
$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,2).Value = "2025-10-14"
$ws.Cells.Item(5,2).ClearFormats()
$ws.Cells.Item(5,3).Value = "19:30:00"
$ws.Cells.Item(5,4).Value = "Chapecoense"
$ws.Cells.Item(5,5).Value = "Botafogo SP"
$ws.Cells.Item(5,6).Value = 1.65
$ws.Cells.Item(5,7).Value = 1.72
$ws.Cells.Item(5,8).Value = 6.2
$ws.Cells.Item(5,9).Value = 7.2
$ws.Cells.Item(5,10).Value = 3.75
$ws.Cells.Item(5,11).Value = 4.2
$ws.Cells.Item(5,12).Value = 1.47
$ws.Cells.Item(5,13).Value = 1.08
$ws.Cells.Item(5,14).Value = 3.1
$ws.Cells.Item(5,15).Value = 1.41
$ws.Cells.Item(5,16).Value = 1.73
$ws.Cells.Item(5,17).Value = 2.18
$ws.Cells.Item(5,18).Value = 1.27
$ws.Cells.Item(5,19).Value = 4.1
$ws.Cells.Item(5,20).Value = 2.12
$ws.Cells.Item(5,21).Value = 1.78
$ws.Cells.Item(5,22).Value = 1.16
$ws.Cells.Item(5,23).Value = 2.38
$ws.Cells.Item(5,24).Value = 14.5
$ws.Cells.Item(5,25).Value = 21
$ws.Cells.Item(5,26).Value = 1000
$ws.Cells.Item(5,27).Value = 260
$ws.Cells.Item(5,28).Value = 7.8
$ws.Cells.Item(5,29).Value = 10.5
$ws.Cells.Item(5,30).Value = 32
$ws.Cells.Item(5,31).Value = 140
$ws.Cells.Item(5,32).Value = 9.800000000000001
$ws.Cells.Item(5,33).Value = 12
$ws.Cells.Item(5,34).Value = 30
$ws.Cells.Item(5,35).Value = 140
$ws.Cells.Item(5,36).Value = 19.5
$ws.Cells.Item(5,37).Value = 24
$ws.Cells.Item(5,38).Value = 60
$ws.Cells.Item(5,39).Value = 230
$ws.Cells.Item(5,40).Value = 15.5
$ws.Cells.Item(5,41).Value = 1000

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = "FIFA World Cup Qualifiers - Americas"
$ws.Cells.Item(6,2).NumberFormat = "@"
$ws.Cells.Item(6,2).Value = "2025-10-14"
$ws.Cells.Item(6,2).ClearFormats()
$ws.Cells.Item(6,3).Value = "20:00:00"
$ws.Cells.Item(6,4).Value = "Curacao"
$ws.Cells.Item(6,5).Value = "Trinidad & Tobago"
$ws.Cells.Item(6,6).Value = 1.71
$ws.Cells.Item(6,7).Value = 1.77
$ws.Cells.Item(6,8).Value = 5.3
$ws.Cells.Item(6,9).Value = 6.8
$ws.Cells.Item(6,10).Value = 3.6
$ws.Cells.Item(6,11).Value = 4.3
$ws.Cells.Item(6,12).Value = 1.37
$ws.Cells.Item(6,13).Value = 1.08
$ws.Cells.Item(6,14).Value = 2.84
$ws.Cells.Item(6,15).Value = 1.36
$ws.Cells.Item(6,16).Value = 1.64
$ws.Cells.Item(6,17).Value = 1.92
$ws.Cells.Item(6,18).Value = 1.28
$ws.Cells.Item(6,19).Value = 3.4
$ws.Cells.Item(6,20).Value = 1.94
$ws.Cells.Item(6,21).Value = 1.83
$ws.Cells.Item(6,22).Value = 1.17
$ws.Cells.Item(6,23).Value = 2.28
$ws.Cells.Item(6,24).Value = 1000
$ws.Cells.Item(6,25).Value = 1000
$ws.Cells.Item(6,26).Value = 1000
$ws.Cells.Item(6,27).Value = 1000
$ws.Cells.Item(6,28).Value = 1000
$ws.Cells.Item(6,29).Value = 1000
$ws.Cells.Item(6,30).Value = 29
$ws.Cells.Item(6,31).Value = 1000
$ws.Cells.Item(6,32).Value = 1000
$ws.Cells.Item(6,33).Value = 1000
$ws.Cells.Item(6,34).Value = 1000
$ws.Cells.Item(6,35).Value = 120
$ws.Cells.Item(6,36).Value = 1000
$ws.Cells.Item(6,37).Value = 1000
$ws.Cells.Item(6,38).Value = 1000
$ws.Cells.Item(6,39).Value = 1000
$ws.Cells.Item(6,40).Value = 1000
$ws.Cells.Item(6,41).Value = 1000

Write-Output "Edit complete"
